$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = 0.26
$ws.Range("G6").Value = 0.25
$ws.Range("J6").Value = 0.18
$ws.Range("M6").Value = 0.21
$ws.Range("P6").Value = 0.24
